$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# Price cells are forced to Text format before assignment so that values such as
# '1.00' or '28.91' are not silently reinterpreted by Excel as numbers, then the
# cell style is restored to Normal so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.813.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.679.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.16%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.42%  '

$ws.Range("E6").Value = '  +3.32%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.91'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.264'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0644'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.01%  '

$ws.Range("E11").Value = '  +0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.921.85'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.679.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.603'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.39%  '

$ws.Range("E16").Value = '  +5.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.816.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.82%  '

$ws.Range("E23").Value = '  +1.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.55%  '

$ws.Range("E27").Value = '  +2.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.08%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("E31").Value = '  +3.73%  '

$ws.Range("E32").Value = '  +3.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.514.87'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.39%  '

$ws.Range("E34").Value = '  +4.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.75'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.81%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '84.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.72%  '

$ws.Range("E37").Value = '  +1.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.603'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.56%  '

$ws.Range("E39").Value = '  +5.04%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  -3.74%  '

$ws.Range("E42").Value = '  +2.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0501'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.47%  '

$ws.Range("E45").Value = '  +2.01%  '

$ws.Range("E46").Value = '  -0.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '50.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.812.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0118'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '92.79'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.17%  '

